$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5513548631581955
$ws.Range("C2").Value = 0.1624190934874861
$ws.Range("D2").Value = 0.2082656351497434
$ws.Range("F2").Value = 1.372022138432399
$ws.Range("G2").Value = 0.7516564130985017
$ws.Range("H2").Value = 0.8580028273867129
$ws.Range("I2").Value = 0.885749898241535
$ws.Range("J2").Value = 0.2057898918912713
$ws.Range("K2").Value = 0.3242496028636594
$ws.Range("L2").Value = 0.3476246894521893
$ws.Range("O2").Value = 3.228388186044896

$ws.Range("B3").Value = 0.5071859174208839
$ws.Range("C3").Value = 0.1618322943610764
$ws.Range("D3").Value = 0.2054934232637038
$ws.Range("F3").Value = 1.380446964905893
$ws.Range("G3").Value = 0.7581792986712941
$ws.Range("H3").Value = 0.8646797307768637
$ws.Range("I3").Value = 0.8945297417896505
$ws.Range("J3").Value = 0.2071058191985919
$ws.Range("K3").Value = 0.2854423381356241
$ws.Range("L3").Value = 0.3416159085320345
$ws.Range("O3").Value = 3.255956297056727

$ws.Range("B4").Value = 0.4801105707411466
$ws.Range("C4").Value = 0.1614883107873837
$ws.Range("D4").Value = 0.2038638569240305
$ws.Range("F4").Value = 1.386259559993611
$ws.Range("G4").Value = 0.7626067025994772
$ws.Range("H4").Value = 0.8690971644646552
$ws.Range("I4").Value = 0.9003183030509945
$ws.Range("J4").Value = 0.2079874840256437
$ws.Range("K4").Value = 0.2615550756848108
$ws.Range("L4").Value = 0.3380536345962497
$ws.Range("O4").Value = 3.274435550144787

$ws.Range("B5").Value = 0.4690892446665771
$ws.Range("C5").Value = 0.1613522728317065
$ws.Range("D5").Value = 0.2032181533383124
$ws.Range("F5").Value = 1.388789268338662
$ws.Range("G5").Value = 0.7645171223081135
$ws.Range("H5").Value = 0.870977301197378
$ws.Range("I5").Value = 0.9027772466252273
$ws.Range("J5").Value = 0.2083653249548796
$ws.Range("K5").Value = 0.2518065839901311
$ws.Range("L5").Value = 0.3366341325845781
$ws.Range("O5").Value = 3.282356555811461

$ws.Range("B6").Value = 0.4672599187930473
$ws.Range("C6").Value = 0.1613299348286077
$ws.Range("D6").Value = 0.2031120466508298
$ws.Range("F6").Value = 1.389219055103766
$ws.Range("G6").Value = 0.7648407618370925
$ws.Range("H6").Value = 0.871294330606041
$ws.Range("I6").Value = 0.9031915971237083
$ws.Range("J6").Value = 0.2084291866700454
$ws.Range("K6").Value = 0.2501870153295442
$ws.Range("L6").Value = 0.3364003727245688
$ws.Range("O6").Value = 3.283695427978799

$ws.Range("B7").Value = 0.4799618831181078
$ws.Range("C7").Value = 0.1614864593302272
$ws.Range("D7").Value = 0.2038550742675511
$ws.Range("F7").Value = 1.386293024304784
$ws.Range("G7").Value = 0.7626320370946118
$ws.Range("H7").Value = 0.869122196628112
$ws.Range("I7").Value = 0.900351060023489
$ws.Range("J7").Value = 0.2079925045536299
$ws.Range("K7").Value = 0.2614236608250451
$ws.Range("L7").Value = 0.3380343602643308
$ws.Range("O7").Value = 3.274540793947253

$ws.Range("B8").Value = 0.536116681562163
$ws.Range("C8").Value = 0.162213400534192
$ws.Range("D8").Value = 0.2072947744558888
$ws.Range("F8").Value = 1.374794344372184
$ws.Range("G8").Value = 0.7538178618627782
$ws.Range("H8").Value = 0.8602391187733289
$ws.Range("I8").Value = 0.8886946773240929
$ws.Range("J8").Value = 0.2062283480359426
$ws.Range("K8").Value = 0.3108816092639017
$ws.Range("L8").Value = 0.3455265805878582
$ws.Range("O8").Value = 3.2375716163709

$ws.Range("B9").Value = 0.6465559809495005
$ws.Range("C9").Value = 0.1637669042175958
$ws.Range("D9").Value = 0.2146117495458526
$ws.Range("F9").Value = 1.357314480962657
$ws.Range("G9").Value = 0.7398836024855342
$ws.Range("H9").Value = 0.8453370842202048
$ws.Range("I9").Value = 0.86898969753668
$ws.Range("J9").Value = 0.2033522471731359
$ws.Range("K9").Value = 0.4073702464062308
$ws.Range("L9").Value = 0.3612205951800149
$ws.Range("O9").Value = 3.177382950138593

$ws.Range("B10").Value = 0.7278535343180295
$ws.Range("C10").Value = 0.1649844788054153
$ws.Range("D10").Value = 0.2203311820310034
$ws.Range("F10").Value = 1.347553683468803
$ws.Range("G10").Value = 0.7316879574948771
$ws.Range("H10").Value = 0.8359181206427095
$ws.Range("I10").Value = 0.8564312346105325
$ws.Range("J10").Value = 0.2015932186799105
$ws.Range("K10").Value = 0.4779277429140905
$ws.Range("L10").Value = 0.3733537536396909
$ws.Range("O10").Value = 3.140653576254991

$ws.Range("B11").Value = 0.7648647717423103
$ws.Range("C11").Value = 0.1655545665190132
$ws.Range("D11").Value = 0.2230067263702011
$ws.Range("F11").Value = 1.343780635219616
$ws.Range("G11").Value = 0.7284028159038058
$ws.Range("H11").Value = 0.8319642104098719
$ws.Range("I11").Value = 0.8511339788050307
$ws.Range("J11").Value = 0.2008695286905073
$ws.Range("K11").Value = 0.5099482162065954
$ws.Range("L11").Value = 0.3790027283535835
$ws.Range("O11").Value = 3.125568625222158

$ws.Range("B12").Value = 0.778883274749802
$ws.Range("C12").Value = 0.1657727408642771
$ws.Range("D12").Value = 0.2240303912713557
$ws.Range("F12").Value = 1.342447668119625
$ws.Range("G12").Value = 0.7272225258208351
$ws.Range("H12").Value = 0.8305144557999924
$ws.Range("I12").Value = 0.8491877667926104
$ws.Range("J12").Value = 0.200606459650448
$ws.Range("K12").Value = 0.5220619072255772
$ws.Range("L12").Value = 0.3811603097088465
$ws.Range("O12").Value = 3.120089603774503

$ws.Range("B13").Value = 0.7758640144096489
$ws.Range("C13").Value = 0.1657256516378141
$ws.Range("D13").Value = 0.2238094617371189
$ws.Range("F13").Value = 1.342730487686566
$ws.Range("G13").Value = 0.7274738882992153
$ws.Range("H13").Value = 0.8308245746383704
$ws.Range("I13").Value = 0.8496042616720025
$ws.Range("J13").Value = 0.200662628468411
$ws.Range("K13").Value = 0.5194535426934124
$ws.Range("L13").Value = 0.38069481899322
$ws.Range("O13").Value = 3.121259234004285

$ws.Range("B14").Value = 0.7660180247583526
$ws.Range("C14").Value = 0.1655724700736414
$ws.Range("D14").Value = 0.2230907342076307
$ws.Range("F14").Value = 1.343669051934533
$ws.Range("G14").Value = 0.7283044356249988
$ws.Range("H14").Value = 0.8318439865277796
$ws.Range("I14").Value = 0.8509726655682925
$ws.Range("J14").Value = 0.2008476659913576
$ws.Range("K14").Value = 0.5109450570872411
$ws.Range("L14").Value = 0.3791798654518601
$ws.Range("O14").Value = 3.125113187423338

$ws.Range("B15").Value = 0.7599874537105222
$ws.Range("C15").Value = 0.165478939690189
$ws.Range("D15").Value = 0.2226518562427202
$ws.Range("F15").Value = 1.344256422139331
$ws.Range("G15").Value = 0.7288214684240373
$ws.Range("H15").Value = 0.8324745903753694
$ws.Range("I15").Value = 0.8518186320871806
$ws.Range("J15").Value = 0.2009624355877264
$ws.Range("K15").Value = 0.5057318087491467
$ws.Range("L15").Value = 0.3782543080820773
$ws.Range("O15").Value = 3.127504226891659

$ws.Range("B16").Value = 0.7254352559147605
$ws.Range("C16").Value = 0.1649475452712892
$ws.Range("D16").Value = 0.2201578034477762
$ws.Range("F16").Value = 1.347813674919486
$ws.Range("G16").Value = 0.7319115659167821
$ws.Range("H16").Value = 0.8361831692876649
$ws.Range("I16").Value = 0.8567857844548605
$ws.Range("J16").Value = 0.2016420508566519
$ws.Range("K16").Value = 0.4758335248287437
$ws.Range("L16").Value = 0.3729871702541061
$ws.Range("O16").Value = 3.14167207013098

$ws.Range("B17").Value = 0.7042452374560924
$ws.Range("C17").Value = 0.1646256772663435
$ws.Range("D17").Value = 0.2186465924790326
$ws.Range("F17").Value = 1.350166720153638
$ws.Range("G17").Value = 0.7339207323572481
$ws.Range("H17").Value = 0.8385429466676371
$ws.Range("I17").Value = 0.859939411168881
$ws.Range("J17").Value = 0.2020785493756172
$ws.Range("K17").Value = 0.4574717585765313
$ws.Range("L17").Value = 0.3697889946137849
$ws.Range("O17").Value = 3.150779292963094

$ws.Range("B18").Value = 0.6920600458111323
$ws.Range("C18").Value = 0.1644420749125501
$ws.Range("D18").Value = 0.2177843317948742
$ws.Range("F18").Value = 1.351582936569052
$ws.Range("G18").Value = 0.7351180580349137
$ws.Range("H18").Value = 0.839931370845143
$ws.Range("I18").Value = 0.861792418557183
$ws.Range("J18").Value = 0.2023368138420878
$ws.Range("K18").Value = 0.4469034028908538
$ws.Range("L18").Value = 0.3679616928921803
$ws.Range("O18").Value = 3.156170338634666

$ws.Range("B19").Value = 0.6879348502371556
$ws.Range("C19").Value = 0.1643801737805362
$ws.Range("D19").Value = 0.2174935818096344
$ws.Range("F19").Value = 1.352073234428111
$ws.Range("G19").Value = 0.7355306148756711
$ws.Range("H19").Value = 0.8404068182298374
$ws.Range("I19").Value = 0.8624265355332454
$ws.Range("J19").Value = 0.2024254954489386
$ws.Range("K19").Value = 0.4433239341024091
$ws.Range("L19").Value = 0.3673451018317877
$ws.Range("O19").Value = 3.158021903452664

$ws.Range("B20").Value = 0.7065006737877866
$ws.Range("C20").Value = 0.1646597828246641
$ws.Range("D20").Value = 0.2188067451844233
$ws.Range("F20").Value = 1.34990973516566
$ws.Range("G20").Value = 0.7337025366797292
$ws.Range("H20").Value = 0.8382885216739169
$ws.Range("I20").Value = 0.859599653152145
$ws.Range("J20").Value = 0.2020313381512437
$ws.Range("K20").Value = 0.4594271458636001
$ws.Range("L20").Value = 0.3701281841494506
$ws.Range("O20").Value = 3.149793998977785

$ws.Range("B21").Value = 0.7689099527068208
$ws.Range("C21").Value = 0.1656174012262568
$ws.Range("D21").Value = 0.223301557970629
$ws.Range("F21").Value = 1.343390773849926
$ws.Range("G21").Value = 0.7280587543187877
$ws.Range("H21").Value = 0.8315432717920643
$ws.Range("I21").Value = 0.8505691109548437
$ws.Range("J21").Value = 0.2007930182743749
$ws.Range("K21").Value = 0.5134445310351339
$ws.Range("L21").Value = 0.3796243451674002
$ws.Range("O21").Value = 3.123974856584312

$ws.Range("B22").Value = 0.8097159455457472
$ws.Range("C22").Value = 0.1662566177815847
$ws.Range("D22").Value = 0.2263002957725462
$ws.Range("F22").Value = 1.339688613378222
$ws.Range("G22").Value = 0.724741637060113
$ws.Range("H22").Value = 0.8274117280802784
$ws.Range("I22").Value = 0.8450153498671291
$ws.Range("J22").Value = 0.2000476738389381
$ws.Range("K22").Value = 0.5486790708943659
$ws.Range("L22").Value = 0.385938013184898
$ws.Range("O22").Value = 3.10846047595399

$ws.Range("B23").Value = 0.7879356953819467
$ws.Range("C23").Value = 0.1659142454947897
$ws.Range("D23").Value = 0.2246942571088084
$ws.Range("F23").Value = 1.341613481540485
$ws.Range("G23").Value = 0.7264780585617174
$ws.Range("H23").Value = 0.8295914993074405
$ws.Range("I23").Value = 0.8479476432460693
$ws.Range("J23").Value = 0.2004396329511273
$ws.Range("K23").Value = 0.5298802956901625
$ws.Range("L23").Value = 0.3825585276147478
$ws.Range("O23").Value = 3.116616401331555

$ws.Range("B24").Value = 0.7054809995712787
$ws.Range("C24").Value = 0.1646443592029456
$ws.Range("D24").Value = 0.2187343197213778
$ws.Range("F24").Value = 1.350025720523249
$ws.Range("G24").Value = 0.7338010514058979
$ws.Range("H24").Value = 0.8384034482889646
$ws.Range("I24").Value = 0.8597531333231174
$ws.Range("J24").Value = 0.2020526595565819
$ws.Range("K24").Value = 0.4585431523786951
$ws.Range("L24").Value = 0.3699748011243997
$ws.Range("O24").Value = 3.150238967000249

$ws.Range("B25").Value = 0.6166488331266464
$ws.Range("C25").Value = 0.1633331269967613
$ws.Range("D25").Value = 0.2125716118413408
$ws.Range("F25").Value = 1.361501456691272
$ws.Range("G25").Value = 0.7432945733900667
$ws.Range("H25").Value = 0.8490994931517974
$ws.Range("I25").Value = 0.8739832081630396
$ws.Range("J25").Value = 0.2040680184820687
$ws.Range("K25").Value = 0.3813238568168345
$ws.Range("L25").Value = 0.3568685217386758
$ws.Range("O25").Value = 3.19234917665402
